$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns L:O need the same style as the existing header (A1).
# Copy format from A1 (bold font, border, centered alignment) into the new header cells
# before overwriting their values, so they pick up the existing style index
# instead of Excel allocating a brand-new (duplicate) style.
$ws.Range("A1").Copy($ws.Range("L1:O1"))

# Header row (row 1)
$ws.Range("A1").Value = 'Best Estimator'
$ws.Range("B1").Value = 'Best Score'
$ws.Range("C1").Value = 'Best Params'
$ws.Range("D1").Value = 'CV Train F1'
$ws.Range("E1").Value = 'CV Test F1'
$ws.Range("F1").Value = 'Validation F1'
$ws.Range("G1").Value = 'CV Train Precision'
$ws.Range("H1").Value = 'CV Test Precision'
$ws.Range("I1").Value = 'Validation Precision'
$ws.Range("J1").Value = 'CV Train Recall'
$ws.Range("K1").Value = 'CV Test Recall'
$ws.Range("L1").Value = 'Validation Recall'
$ws.Range("M1").Value = 'Y Val (Validation)'
$ws.Range("N1").Value = 'Y Pred (Validation)'
$ws.Range("O1").Value = 'Seed'

# Data rows 2-6
# Row 2
$ws.Range("A2").Value = 'Pipeline(steps=[(''scaler'', RobustScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7fd5873896d0>),
                (''model'',
                 RandomForestClassifier(max_depth=3, min_samples_leaf=5,
                                        min_samples_split=5, n_estimators=5,
                                        random_state=42))])'
$ws.Range("B2").Value = 0.6759523809523811
$ws.Range("C2").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7fd5872e18e0>, ''scaler'': RobustScaler(), ''model__n_estimators'': 5, ''model__min_samples_split'': 5, ''model__min_samples_leaf'': 5, ''model__max_features'': ''sqrt'', ''model__max_depth'': 3, ''model__class_weight'': None}'
$ws.Range("D2").Value = 0.8513278198212739
$ws.Range("E2").Value = 0.5244346542346543
$ws.Range("F2").Value = 0.7692307692307692
$ws.Range("G2").Value = 0.8235022318487288
$ws.Range("H2").Value = 0.5440571428571428
$ws.Range("I2").Value = 0.6666666666666666
$ws.Range("J2").Value = 0.8901190476190477
$ws.Range("K2").Value = 0.5496
$ws.Range("L2").Value = 0.9090909090909091
$ws.Range("M2").Value = '[1 0 1 1 1 1 0 1 0 1 0 1 0 1 1 0 0 1 1 1 1 0 1 1 0 1 1 1 1 0 0 0 0 1 0 1]'
$ws.Range("N2").Value = '[1 1 1 1 1 0 0 1 1 1 1 1 1 1 0 0 0 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 0 1 1 1]'
$ws.Range("O2").Value = 42

# Row 3
$ws.Range("A3").Value = 'Pipeline(steps=[(''scaler'', RobustScaler()), (''selector'', None),
                (''model'',
                 RandomForestClassifier(max_depth=6, max_features=''log2'',
                                        min_samples_leaf=2, n_estimators=5,
                                        random_state=42))])'
$ws.Range("B3").Value = 0.6921428571428571
$ws.Range("C3").Value = '{''selector'': None, ''scaler'': RobustScaler(), ''model__n_estimators'': 5, ''model__min_samples_split'': 2, ''model__min_samples_leaf'': 2, ''model__max_features'': ''log2'', ''model__max_depth'': 6, ''model__class_weight'': None}'
$ws.Range("D3").Value = 0.839690148833381
$ws.Range("E3").Value = 0.5901477744477744
$ws.Range("F3").Value = 0.5416666666666666
$ws.Range("G3").Value = 0.8243197915927847
$ws.Range("H3").Value = 0.5678964285714286
$ws.Range("I3").Value = 0.5416666666666666
$ws.Range("J3").Value = 0.8633749999999999
$ws.Range("K3").Value = 0.6416
$ws.Range("L3").Value = 0.5416666666666666
$ws.Range("M3").Value = '[1 1 0 1 0 0 1 0 1 1 1 0 1 1 1 1 1 1 1 1 0 0 1 0 0 1 0 1 1 0 1 1 0 1 1 1]'
$ws.Range("N3").Value = '[0 1 1 1 1 0 0 1 0 0 0 1 1 0 1 1 1 1 0 0 1 1 1 1 1 1 1 0 1 1 1 0 1 0 1 1]'
$ws.Range("O3").Value = 69

# Row 4
$ws.Range("A4").Value = 'Pipeline(steps=[(''scaler'', RobustScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7fd5872e1eb0>),
                (''model'',
                 RandomForestClassifier(max_depth=2, min_samples_leaf=2,
                                        min_samples_split=4, n_estimators=5,
                                        random_state=42))])'
$ws.Range("B4").Value = 0.6713095238095238
$ws.Range("C4").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7fd58720c460>, ''scaler'': RobustScaler(), ''model__n_estimators'': 5, ''model__min_samples_split'': 4, ''model__min_samples_leaf'': 2, ''model__max_features'': ''sqrt'', ''model__max_depth'': 2, ''model__class_weight'': None}'
$ws.Range("D4").Value = 0.8615067929819722
$ws.Range("E4").Value = 0.5544539405039405
$ws.Range("F4").Value = 0.64
$ws.Range("G4").Value = 0.8621157640155077
$ws.Range("H4").Value = 0.5555412698412698
$ws.Range("I4").Value = 0.6666666666666666
$ws.Range("J4").Value = 0.8652368421052631
$ws.Range("K4").Value = 0.5803999999999999
$ws.Range("L4").Value = 0.6153846153846154
$ws.Range("M4").Value = '[0 1 0 0 1 1 1 1 1 1 1 0 1 1 1 1 1 1 1 1 0 1 1 1 0 1 0 1 0 1 0 1 1 1 0 1]'
$ws.Range("N4").Value = '[1 1 1 1 0 1 0 1 1 0 0 0 0 1 1 0 0 1 1 0 1 0 1 1 0 0 1 1 1 1 1 1 1 1 1 1]'
$ws.Range("O4").Value = 23

# Row 5
$ws.Range("A5").Value = 'Pipeline(steps=[(''scaler'', RobustScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7fd58720cdf0>),
                (''model'',
                 RandomForestClassifier(class_weight=''balanced'', max_depth=6,
                                        min_samples_leaf=6, n_estimators=50,
                                        random_state=42))])'
$ws.Range("B5").Value = 0.7155952380952381
$ws.Range("C5").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7fd5870d20d0>, ''scaler'': RobustScaler(), ''model__n_estimators'': 50, ''model__min_samples_split'': 2, ''model__min_samples_leaf'': 6, ''model__max_features'': ''sqrt'', ''model__max_depth'': 6, ''model__class_weight'': ''balanced''}'
$ws.Range("D5").Value = 0.8473978917229471
$ws.Range("E5").Value = 0.5674308302808304
$ws.Range("F5").Value = 0.6363636363636364
$ws.Range("G5").Value = 0.8223892449520005
$ws.Range("H5").Value = 0.542618253968254
$ws.Range("I5").Value = 0.6363636363636364
$ws.Range("J5").Value = 0.8873095238095238
$ws.Range("K5").Value = 0.6314
$ws.Range("L5").Value = 0.6363636363636364
$ws.Range("M5").Value = '[0 1 1 0 0 1 0 0 0 0 1 1 1 0 0 1 1 0 1 1 1 1 1 1 1 1 0 0 1 0 1 1 1 1 1 0]'
$ws.Range("N5").Value = '[0 1 1 1 0 0 0 1 0 0 1 1 0 0 1 0 0 1 1 1 0 1 1 1 1 0 1 1 1 1 0 1 1 1 0 1]'
$ws.Range("O5").Value = 99

# Row 6
$ws.Range("A6").Value = 'Pipeline(steps=[(''scaler'', StandardScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7fd58722a970>),
                (''model'',
                 RandomForestClassifier(class_weight=''balanced'', max_depth=1,
                                        min_samples_leaf=3, min_samples_split=3,
                                        n_estimators=5, random_state=42))])'
$ws.Range("B6").Value = 0.7030952380952381
$ws.Range("C6").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7fd5870e16a0>, ''scaler'': StandardScaler(), ''model__n_estimators'': 5, ''model__min_samples_split'': 3, ''model__min_samples_leaf'': 3, ''model__max_features'': ''sqrt'', ''model__max_depth'': 1, ''model__class_weight'': ''balanced''}'
$ws.Range("D6").Value = 0.8543515421504925
$ws.Range("E6").Value = 0.6273597791097791
$ws.Range("F6").Value = 0.55
$ws.Range("G6").Value = 0.8193993821591412
$ws.Range("H6").Value = 0.569670634920635
$ws.Range("I6").Value = 0.55
$ws.Range("J6").Value = 0.9027272727272728
$ws.Range("K6").Value = 0.7260000000000001
$ws.Range("L6").Value = 0.55
$ws.Range("M6").Value = '[1 0 1 1 0 0 0 0 1 0 1 1 0 1 1 0 1 0 0 0 0 0 1 1 1 0 1 0 1 1 1 1 1 1 1 0]'
$ws.Range("N6").Value = '[1 1 1 0 1 1 0 0 1 0 0 0 1 1 1 0 1 1 0 1 1 1 1 1 0 0 1 1 0 1 0 0 1 0 0 0]'
$ws.Range("O6").Value = 89

